# Lithuania A Lyga bases update (11-04-2024 23:56):
#  - rows 26/27, 100/102 and 103/104 had their match data (columns B..AC)
#    swapped with each other (column A, the running index, stays put);
#  - four brand-new fixtures were appended as rows 136-139.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 26 and row 27 (all columns except A) ---
$ws.Cells.Item(26, 2).Value = 6732773
$ws.Cells.Item(26, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(26, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(26, 5).Value = 45109.58333333334
$ws.Cells.Item(26, 6).Value = "Suduva Marijampole"
$ws.Cells.Item(26, 7).Value = "Hegelmann Litauen"
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = "A"
$ws.Cells.Item(26, 11).Value = 5
$ws.Cells.Item(26, 12).Value = 3.8
$ws.Cells.Item(26, 13).Value = 1.533
$ws.Cells.Item(26, 14).Value = 5
$ws.Cells.Item(26, 15).Value = 4.2
$ws.Cells.Item(26, 16).Value = 1.533
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = 1.875
$ws.Cells.Item(26, 19).Value = 1.925
$ws.Cells.Item(26, 20).Value = 2.5
$ws.Cells.Item(26, 21).Value = 1.9
$ws.Cells.Item(26, 22).Value = 1.9
$ws.Cells.Item(26, 23).Value = -1
$ws.Cells.Item(26, 24).Value = -1
$ws.Cells.Item(26, 25).Value = 0.5329999999999999
$ws.Cells.Item(26, 26).Value = 0
$ws.Cells.Item(26, 27).Value = -0
$ws.Cells.Item(26, 28).Value = -1
$ws.Cells.Item(26, 29).Value = 0.8999999999999999

$ws.Cells.Item(27, 2).Value = 6732711
$ws.Cells.Item(27, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(27, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(27, 5).Value = 45109.58333333334
$ws.Cells.Item(27, 6).Value = "Banga Gargzdai"
$ws.Cells.Item(27, 7).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(27, 8).Value = 1
$ws.Cells.Item(27, 9).Value = 4
$ws.Cells.Item(27, 10).Value = "A"
$ws.Cells.Item(27, 11).Value = 5
$ws.Cells.Item(27, 12).Value = 3.6
$ws.Cells.Item(27, 13).Value = 1.571
$ws.Cells.Item(27, 14).Value = 11
$ws.Cells.Item(27, 15).Value = 4.75
$ws.Cells.Item(27, 16).Value = 1.25
$ws.Cells.Item(27, 17).Value = 1.5
$ws.Cells.Item(27, 18).Value = 1.975
$ws.Cells.Item(27, 19).Value = 1.825
$ws.Cells.Item(27, 20).Value = 2.5
$ws.Cells.Item(27, 21).Value = 1.8
$ws.Cells.Item(27, 22).Value = 2
$ws.Cells.Item(27, 23).Value = -1
$ws.Cells.Item(27, 24).Value = -1
$ws.Cells.Item(27, 25).Value = 0.25
$ws.Cells.Item(27, 26).Value = -1
$ws.Cells.Item(27, 27).Value = 0.825
$ws.Cells.Item(27, 28).Value = 0.8
$ws.Cells.Item(27, 29).Value = -1

# --- Swap row 100 and row 102 (all columns except A) ---
$ws.Cells.Item(100, 2).Value = 6732727
$ws.Cells.Item(100, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(100, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(100, 5).Value = 45242.41319444445
$ws.Cells.Item(100, 6).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(100, 7).Value = "FK Dainava Alytus"
$ws.Cells.Item(100, 8).Value = 1
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = "H"
$ws.Cells.Item(100, 11).Value = 1.285
$ws.Cells.Item(100, 12).Value = 5.5
$ws.Cells.Item(100, 13).Value = 6.5
$ws.Cells.Item(100, 14).Value = 1.3
$ws.Cells.Item(100, 15).Value = 5.5
$ws.Cells.Item(100, 16).Value = 6
$ws.Cells.Item(100, 17).Value = -1.5
$ws.Cells.Item(100, 18).Value = 1.9
$ws.Cells.Item(100, 19).Value = 1.9
$ws.Cells.Item(100, 20).Value = 2.75
$ws.Cells.Item(100, 21).Value = 1.8
$ws.Cells.Item(100, 22).Value = 2
$ws.Cells.Item(100, 23).Value = 0.3
$ws.Cells.Item(100, 24).Value = -1
$ws.Cells.Item(100, 25).Value = -1
$ws.Cells.Item(100, 26).Value = -1
$ws.Cells.Item(100, 27).Value = 0.8999999999999999
$ws.Cells.Item(100, 28).Value = -1
$ws.Cells.Item(100, 29).Value = 1

$ws.Cells.Item(102, 2).Value = 6732834
$ws.Cells.Item(102, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(102, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(102, 5).Value = 45242.41319444445
$ws.Cells.Item(102, 6).Value = "Panevezys"
$ws.Cells.Item(102, 7).Value = "FK Dziugas Telsiai"
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = "D"
$ws.Cells.Item(102, 11).Value = 1.25
$ws.Cells.Item(102, 12).Value = 5.5
$ws.Cells.Item(102, 13).Value = 7.5
$ws.Cells.Item(102, 14).Value = 1.45
$ws.Cells.Item(102, 15).Value = 4.5
$ws.Cells.Item(102, 16).Value = 5
$ws.Cells.Item(102, 17).Value = -1
$ws.Cells.Item(102, 18).Value = 1.775
$ws.Cells.Item(102, 19).Value = 2.025
$ws.Cells.Item(102, 20).Value = 2.5
$ws.Cells.Item(102, 21).Value = 1.875
$ws.Cells.Item(102, 22).Value = 1.925
$ws.Cells.Item(102, 23).Value = -1
$ws.Cells.Item(102, 24).Value = 3.5
$ws.Cells.Item(102, 25).Value = -1
$ws.Cells.Item(102, 26).Value = -1
$ws.Cells.Item(102, 27).Value = 1.025
$ws.Cells.Item(102, 28).Value = -1
$ws.Cells.Item(102, 29).Value = 0.925

# --- Swap row 103 and row 104 (all columns except A) ---
$ws.Cells.Item(103, 2).Value = 6732836
$ws.Cells.Item(103, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(103, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(103, 5).Value = 45242.41319444445
$ws.Cells.Item(103, 6).Value = "FK Siauliai"
$ws.Cells.Item(103, 7).Value = "Banga Gargzdai"
$ws.Cells.Item(103, 8).Value = 3
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = "H"
$ws.Cells.Item(103, 11).Value = 1.222
$ws.Cells.Item(103, 12).Value = 5.5
$ws.Cells.Item(103, 13).Value = 9
$ws.Cells.Item(103, 14).Value = 1.363
$ws.Cells.Item(103, 15).Value = 4.5
$ws.Cells.Item(103, 16).Value = 7
$ws.Cells.Item(103, 17).Value = -1.25
$ws.Cells.Item(103, 18).Value = 1.9
$ws.Cells.Item(103, 19).Value = 1.9
$ws.Cells.Item(103, 20).Value = 2.5
$ws.Cells.Item(103, 21).Value = 1.975
$ws.Cells.Item(103, 22).Value = 1.825
$ws.Cells.Item(103, 23).Value = 0.363
$ws.Cells.Item(103, 24).Value = -1
$ws.Cells.Item(103, 25).Value = -1
$ws.Cells.Item(103, 26).Value = 0.8999999999999999
$ws.Cells.Item(103, 27).Value = -1
$ws.Cells.Item(103, 28).Value = 0.9750000000000001
$ws.Cells.Item(103, 29).Value = -1

$ws.Cells.Item(104, 2).Value = 7465686
$ws.Cells.Item(104, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(104, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(104, 5).Value = 45242.41319444445
$ws.Cells.Item(104, 6).Value = "FK Kauno Zalgiris"
$ws.Cells.Item(104, 7).Value = "Hegelmann Litauen"
$ws.Cells.Item(104, 8).Value = 4
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = "H"
$ws.Cells.Item(104, 11).Value = 2.3
$ws.Cells.Item(104, 12).Value = 4
$ws.Cells.Item(104, 13).Value = 2.3
$ws.Cells.Item(104, 14).Value = 2.55
$ws.Cells.Item(104, 15).Value = 4
$ws.Cells.Item(104, 16).Value = 2.2
$ws.Cells.Item(104, 17).Value = 0.25
$ws.Cells.Item(104, 18).Value = 1.8
$ws.Cells.Item(104, 19).Value = 2
$ws.Cells.Item(104, 20).Value = 2.75
$ws.Cells.Item(104, 21).Value = 1.85
$ws.Cells.Item(104, 22).Value = 1.95
$ws.Cells.Item(104, 23).Value = 1.55
$ws.Cells.Item(104, 24).Value = -1
$ws.Cells.Item(104, 25).Value = -1
$ws.Cells.Item(104, 26).Value = 0.8
$ws.Cells.Item(104, 27).Value = -1
$ws.Cells.Item(104, 28).Value = 0.8500000000000001
$ws.Cells.Item(104, 29).Value = -1

# --- Add new rows 136-139 ---
$ws.Range("A135:AC135").Copy()
$ws.Range("A136:AC139").PasteSpecial(-4122)

$ws.Cells.Item(136, 1).Value = 134
$ws.Cells.Item(136, 2).Value = 7862922
$ws.Cells.Item(136, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(136, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(136, 5).Value = 45392.5
$ws.Cells.Item(136, 6).Value = "FK Siauliai"
$ws.Cells.Item(136, 7).Value = "Panevezys"
$ws.Cells.Item(136, 8).Value = 1
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 10).Value = "D"
$ws.Cells.Item(136, 11).Value = 2.7
$ws.Cells.Item(136, 12).Value = 3
$ws.Cells.Item(136, 13).Value = 2.5
$ws.Cells.Item(136, 14).Value = 2.9
$ws.Cells.Item(136, 15).Value = 2.9
$ws.Cells.Item(136, 16).Value = 2.375
$ws.Cells.Item(136, 17).Value = 0.25
$ws.Cells.Item(136, 18).Value = 1.75
$ws.Cells.Item(136, 19).Value = 2.05
$ws.Cells.Item(136, 20).Value = 1.75
$ws.Cells.Item(136, 21).Value = 1.775
$ws.Cells.Item(136, 22).Value = 2.025
$ws.Cells.Item(136, 23).Value = -1
$ws.Cells.Item(136, 24).Value = 1.9
$ws.Cells.Item(136, 25).Value = -1
$ws.Cells.Item(136, 26).Value = 0.375
$ws.Cells.Item(136, 27).Value = -0.5
$ws.Cells.Item(136, 28).Value = 0.3875
$ws.Cells.Item(136, 29).Value = -0.5

$ws.Cells.Item(137, 1).Value = 135
$ws.Cells.Item(137, 2).Value = 7862044
$ws.Cells.Item(137, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(137, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(137, 5).Value = 45392.5
$ws.Cells.Item(137, 6).Value = "Banga Gargzdai"
$ws.Cells.Item(137, 7).Value = "Suduva Marijampole"
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = "D"
$ws.Cells.Item(137, 11).Value = 2.1
$ws.Cells.Item(137, 12).Value = 3.05
$ws.Cells.Item(137, 13).Value = 3.2
$ws.Cells.Item(137, 14).Value = 3.4
$ws.Cells.Item(137, 15).Value = 3
$ws.Cells.Item(137, 16).Value = 2.1
$ws.Cells.Item(137, 17).Value = 0.25
$ws.Cells.Item(137, 18).Value = 1.95
$ws.Cells.Item(137, 19).Value = 1.85
$ws.Cells.Item(137, 20).Value = 2
$ws.Cells.Item(137, 21).Value = 1.95
$ws.Cells.Item(137, 22).Value = 1.85
$ws.Cells.Item(137, 23).Value = -1
$ws.Cells.Item(137, 24).Value = 2
$ws.Cells.Item(137, 25).Value = -1
$ws.Cells.Item(137, 26).Value = 0.475
$ws.Cells.Item(137, 27).Value = -0.5
$ws.Cells.Item(137, 28).Value = -1
$ws.Cells.Item(137, 29).Value = 0.8500000000000001

$ws.Cells.Item(138, 1).Value = 136
$ws.Cells.Item(138, 2).Value = 7862923
$ws.Cells.Item(138, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(138, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(138, 5).Value = 45392.54166666666
$ws.Cells.Item(138, 6).Value = "FK Kauno Zalgiris"
$ws.Cells.Item(138, 7).Value = "FK Dziugas Telsiai"
$ws.Cells.Item(138, 8).Value = 3
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = "H"
$ws.Cells.Item(138, 11).Value = 1.333
$ws.Cells.Item(138, 12).Value = 4.333
$ws.Cells.Item(138, 13).Value = 7.5
$ws.Cells.Item(138, 14).Value = 1.75
$ws.Cells.Item(138, 15).Value = 3.2
$ws.Cells.Item(138, 16).Value = 4.333
$ws.Cells.Item(138, 17).Value = -0.5
$ws.Cells.Item(138, 18).Value = 1.8
$ws.Cells.Item(138, 19).Value = 2
$ws.Cells.Item(138, 20).Value = 2
$ws.Cells.Item(138, 21).Value = 1.775
$ws.Cells.Item(138, 22).Value = 2.025
$ws.Cells.Item(138, 23).Value = 0.75
$ws.Cells.Item(138, 24).Value = -1
$ws.Cells.Item(138, 25).Value = -1
$ws.Cells.Item(138, 26).Value = 0.8
$ws.Cells.Item(138, 27).Value = -1
$ws.Cells.Item(138, 28).Value = 0.7749999999999999
$ws.Cells.Item(138, 29).Value = -1

$ws.Cells.Item(139, 1).Value = 137
$ws.Cells.Item(139, 2).Value = 7865008
$ws.Cells.Item(139, 3).Value = "Lithuania A Lyga"
$ws.Cells.Item(139, 4).Value = "Lithuania A Lyga"
$ws.Cells.Item(139, 5).Value = 45392.58333333334
$ws.Cells.Item(139, 6).Value = "FK Zalgiris Vilnius"
$ws.Cells.Item(139, 7).Value = "FK Transinvest"
$ws.Cells.Item(139, 8).Value = 2
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = "H"
$ws.Cells.Item(139, 11).Value = 1.333
$ws.Cells.Item(139, 12).Value = 4.2
$ws.Cells.Item(139, 13).Value = 8
$ws.Cells.Item(139, 14).Value = 1.25
$ws.Cells.Item(139, 15).Value = 5
$ws.Cells.Item(139, 16).Value = 8
$ws.Cells.Item(139, 17).Value = -1.5
$ws.Cells.Item(139, 18).Value = 1.85
$ws.Cells.Item(139, 19).Value = 1.95
$ws.Cells.Item(139, 20).Value = 2.75
$ws.Cells.Item(139, 21).Value = 1.825
$ws.Cells.Item(139, 22).Value = 1.975
$ws.Cells.Item(139, 23).Value = 0.25
$ws.Cells.Item(139, 24).Value = -1
$ws.Cells.Item(139, 25).Value = -1
$ws.Cells.Item(139, 26).Value = 0.8500000000000001
$ws.Cells.Item(139, 27).Value = -1
$ws.Cells.Item(139, 28).Value = -1
$ws.Cells.Item(139, 29).Value = 0.9750000000000001
